$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.187.89"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "1.814.82"
$ws.Range("E3").Value = "  +4.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.83"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4431"
$ws.Range("E7").Value = "  +4.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3697"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.62"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07683"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.05"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.249"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.553"
$ws.Range("E15").Value = "  +5.28%  "
$ws.Range("D16").Value = "1.842.65"
$ws.Range("E16").Value = "  +6.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.63"
$ws.Range("E17").Value = "  +6.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001083"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +9.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.50"
$ws.Range("E21").Value = "  +4.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.188"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").Value = "28.267.52"
$ws.Range("E23").Value = "  +3.09%  "
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.047"
$ws.Range("E25").Value = "  -15.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.74"
$ws.Range("E26").Value = "  +2.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.70"
$ws.Range("E27").Value = "  +3.97%  "
$ws.Range("D28").Value = "2.022.09"
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.318"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.14"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  -5.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.864"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09210"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02352"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2172"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06209"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6560"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.196"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.140"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.388"
$ws.Range("E44").Value = "  -1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.85"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6068"
$ws.Range("E46").Value = "  +3.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.763"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.28"
$ws.Range("E49").Value = "  +5.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.153"
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06981"
$ws.Range("E51").Value = "  +2.17%  "
